$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.864.18"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "1.907.11"
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'313.21"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'0.4996"
$ws.Range("E7").Value = "  +3.88%  "
$ws.Range("D8").Value = "'0.3810"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.07281"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "'0.9100"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'0.07647"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "1.902.19"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'5.496"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "'91.92"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D17").Value = "'0.000008733"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "27.896.61"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'14.61"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'5.181"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'6.577"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "'153.03"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").Value = "'1.879"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'2.220"
$ws.Range("E26").Value = "  +3.66%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "'115.34"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "'4.922"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "'0.09017"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("D31").Value = "'3.187"
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").Value = "'1.231"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").Value = "'0.7742"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "'0.02087"
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("D36").Value = "'2.546"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "'0.5565"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").Value = "'3.024"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "'6.908"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'8.491"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'111.70"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.62"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4833"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'0.9991"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'0.06056"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "'0.9082"
$ws.Range("E51").Value = "  +0.92%  "
